# Deploying to gh-pages from @ codeforIATI/codelists@56214b660bf2557f2d99cb7c193581a1aef67d30
#
# The published SectorGroup.xlsx had columns E (codeforiati:category-name)
# and F (codeforiati:group-code) swapped relative to their header labels.
# This script corrects the data by swapping the contents of columns E and F
# for every row (including the header row) so that column E becomes
# "codeforiati:group-code" and column F becomes "codeforiati:category-name",
# matching the upstream fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$rowCount = $used.Rows.Count
$lastRow = $firstRow + $rowCount - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    $eStr = "$eVal"
    $fStr = "$fVal"

    # Write F's old content into E.
    if ($fStr -match '^-?\d+$') {
        $eCell.Value = "'" + $fStr
        $eCell.Style = "Normal"
    } else {
        $eCell.Value = $fStr
    }

    # Write E's old content into F.
    if ($eStr -match '^-?\d+$') {
        $fCell.Value = "'" + $eStr
        $fCell.Style = "Normal"
    } else {
        $fCell.Value = $eStr
    }
}
